# Applies the "10.2 Clustering of Social-Network Graphs" section edit.
#
# 1) The paragraph mark (pilcrow) run-properties of the last "ListParagraph"
#    bullet (the one ending in "... ایزومورف نیست.") lose their
#    <w:rFonts w:hint="cs"/> hint and <w:rtl/> flag.
# 2) Three new paragraphs (a spacer, the "10.2 ..." B-Titr heading, and the
#    new body paragraph) are appended right after that bullet, before the
#    trailing bookmark's paragraph end / sectPr.
#
# The Word object model has no direct property for the paragraph-mark's
# run-level bidi flag, so the safe way to make both edits atomically (and
# keep every existing run's XML byte-for-byte untouched) is to delete the
# whole target paragraph and re-insert its corrected OOXML together with
# the new paragraphs via Range.InsertXML.

$d = $word.ActiveDocument

# Sanity-check we really are about to edit the intended (last) paragraph:
# the bullet that ends the "10.1" discussion and precedes the section break.
$paraCount = $d.Paragraphs.Count
if ($paraCount -ne 109) {
    throw "Unexpected paragraph count $paraCount (expected 109) - target paragraph index may be wrong."
}

$target = $d.Paragraphs.Item(109)
$targetText = $target.Range.Text
if ($targetText -notlike "*نیست*") {
    throw "Paragraph 109 does not look like the expected target paragraph."
}

$targetRange = $target.Range
$insertAt = $targetRange.Start
[void]$targetRange.Delete()

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003101AF" w:rsidRPr="003101AF" w:rsidRDefault="003101AF" w:rsidP="00F83BE5"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:bidi/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:lastRenderedPageBreak/><w:t>آنچه ما آنرا به عنوان دوگان نامیدیم در واقع یک دوتایی واقعی نیست؛</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ز</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>را</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> استفاده از </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">این روش </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ساخت</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">ن گراف در گراف </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G’</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">لزوماً </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ک</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ا</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>زومورف</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> گراف</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> می دهد. یک نمون</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">ه از </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>گراف</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> را بده</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>د</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> که در آن دوتا</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>یی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> از </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G’</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>از</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>نظر</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> گراف</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ا</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>زومورف</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> است و مثال د</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>گر</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> هم بدهید</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> که در آن </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>گراف</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>G</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>’</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">از نظر </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">G </w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ایزومورف </w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ن</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r w:rsidRPr="003101AF"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ست</w:t></w:r><w:r w:rsidR="00F83BE5"><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:bidi/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="B Titr"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>10.2 خوشه بندی گراف شبکه های اجتماعی</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>جنبه مهم شبکه ها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> اجتماع</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ن</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> است که آنها حاو</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> جوامع موجودات</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> هستند که توسط بس</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ار</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> از لبه ها به هم وصل م</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> شوند. به عنوان مثال ، ا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>نها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> با گروه</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> از دوستان در مدرسه </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> گروهها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>یی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> از محققان علاقمند به همان موضوع مطابقت دارد.</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>در</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ن</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> بخش ، خوشه بند</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> نمودار را به عنوان راه</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> برا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> شناسا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>یی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> جوامع در نظر م</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> گ</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ر</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>م</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>به</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> نظر م</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> رسد که تکن</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ک</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>یی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> که در فصل 7 آموخته ا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>م</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ، معمولاً برا</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> مشکل خوشه بند</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> نمودارها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> شبکه ها</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> اجتماع</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin" w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Nazanin"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> نامناسب است.</w:t></w:r></w:p>
'@

$ins = $d.Range($insertAt, $insertAt)
[void]$ins.InsertXML($xml)
